# "Generate Report for Handback" — refresh the localization-status report:
#   - Overview status for zh-cn/de-de flips from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - zh-cn / de-de "Latest Handback DateTime" cells get fresh timestamps
#   - the stale "handback file is not the latest" Error Detail is cleared
#   - a couple of columns are widened / narrowed to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: Status columns for zh-cn / de-de ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: Status column, refreshed handback datetime, cleared error detail ---
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-30 09:07:29"

$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"

# --- de-de sheet: Status column, refreshed handback datetime, cleared error detail ---
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-30 09:07:36"

$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"

# --- Column width adjustments ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
